# ---------------------------------------------------------------------------
# Add "Variables" data-dictionary sheet to the chronic-kidney-disease workbook
# and rename the original data sheet to "Base Dataset".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and create the new "Variables" sheet after it.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Base Dataset"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Variables"

# ---------------------------------------------------------------------------
# 2. Data dictionary content (7 columns x 26 rows incl. header).
# ---------------------------------------------------------------------------
$headers = @("Variable Name","Role","Type","Demographic","Description","Units","Missing Values")

$data = @(
    @("age","Feature","Integer","Age","","year","yes"),
    @("bp","Feature","Integer","","blood pressure","mm/Hg","yes"),
    @("sg","Feature","Categorical","","specific gravity","","yes"),
    @("al","Feature","Categorical","","albumin","","yes"),
    @("su","Feature","Categorical","","sugar","","yes"),
    @("rbc","Feature","Binary","","red blood cells","","yes"),
    @("pc","Feature","Binary","","pus cell","","yes"),
    @("pcc","Feature","Binary","","pus cell clumps","","yes"),
    @("ba","Feature","Binary","","bacteria","","yes"),
    @("bgr","Feature","Integer","","blood glucose random","mgs/dl","yes"),
    @("bu","Feature","Integer","","blood urea","mgs/dl","yes"),
    @("sc","Feature","Continuous","","serum creatinine","mgs/dl","yes"),
    @("sod","Feature","Integer","","sodium","mEq/L","yes"),
    @("pot","Feature","Continuous","","potassium","mEq/L","yes"),
    @("hemo","Feature","Continuous","","hemoglobin","gms","yes"),
    @("pcv","Feature","Integer","","packed cell volume","","yes"),
    @("wbcc","Feature","Integer","","white blood cell count","cells/cmm","yes"),
    @("rbcc","Feature","Continuous","","red blood cell count","millions/cmm","yes"),
    @("htn","Feature","Binary","","hypertension","","yes"),
    @("dm","Feature","Binary","","diabetes mellitus","","yes"),
    @("cad","Feature","Binary","","coronary artery disease","","yes"),
    @("appet","Feature","Binary","","appetite","","yes"),
    @("pe","Feature","Binary","","pedal edema","","yes"),
    @("ane","Feature","Binary","","anemia","","yes"),
    @("class","Target","Binary","","ckd or not ckd","","no")
)

$cols = @("A","B","C","D","E","F","G")

for ($c = 0; $c -lt 7; $c++) {
    $ws2.Range($cols[$c] + "1").Value = $headers[$c]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt 7; $c++) {
        $val = $row[$c]
        if ($val -ne "") {
            $ws2.Range($cols[$c] + $rowNum).Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Styling - build the "data row" look on A2 first (non-bold Arial, dark
#    grey, wrapped + vertically centred), then derive the bold "header" look
#    from it via copy + bold toggle so both styles share the same two new
#    fonts and the engine doesn't fragment the font table further.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Font.Color = 3158064
$ws2.Range("A2").Font.Name = "Arial"
$ws2.Range("A2").VerticalAlignment = -4108
$ws2.Range("A2").WrapText = $true

$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A3:G26").PasteSpecial(-4122) | Out-Null
$ws2.Range("B2:G2").PasteSpecial(-4122) | Out-Null

$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").HorizontalAlignment = -4131

$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("B1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column widths (best-fit approximation).
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 17.02
$ws2.Columns.Item(2).ColumnWidth = 7.31
$ws2.Columns.Item(3).ColumnWidth = 10.74
$ws2.Columns.Item(4).ColumnWidth = 16.02
$ws2.Columns.Item(5).ColumnWidth = 22.45
$ws2.Columns.Item(6).ColumnWidth = 12.74
$ws2.Columns.Item(7).ColumnWidth = 17.74

# ---------------------------------------------------------------------------
# 5. Freeze top row like the base dataset sheet.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 6. AutoFilter + the hidden _xlnm._FilterDatabase defined name.
# ---------------------------------------------------------------------------
$ws2.Range("A1:G26").AutoFilter() | Out-Null
$ws2.Names.Add("_xlnm._FilterDatabase", "=Variables!`$A`$1:`$G`$26") | Out-Null

# ---------------------------------------------------------------------------
# 7. Selections: Base Dataset -> I11, Variables -> E8 (Variables ends up the
#    active / tab-selected sheet, matching the authored workbook).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("I11").Select() | Out-Null

$ws2.Activate()
$ws2.Range("E8").Select() | Out-Null

Write-Output "done"
